$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Mandeep Singh"

# Insert a new column before column A to make room for "matchNo"
$ws.Range("A:A").Insert()

# Populate the new column
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "42nd"
